$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4149273333333334
$ws.Range("H2").Value = 1.244782
$ws.Range("I2").Value = 0.1353844755004719
$ws.Range("J2").Value = 0.1353844755004719
$ws.Range("M2").Value = 0.4655266666666667
$ws.Range("N2").Value = 1.39658
$ws.Range("O2").Value = 0.2411272749309853
$ws.Range("P2").Value = 0.2411272749309853
$ws.Range("Q2").Value = 0.1931597383955556
$ws.Range("R2").Value = 1.73843764556
$ws.Range("S2").Value = 0.03264488964538954
$ws.Range("T2").Value = 0.03264488964538954

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4149273333333334
$ws.Range("H3").Value = 1.244782
$ws.Range("I3").Value = 0.1353844755004719
$ws.Range("J3").Value = 0.1353844755004719
$ws.Range("M3").Value = 1.465099666666666
$ws.Range("N3").Value = 4.395299
$ws.Range("O3").Value = 0.7588727250690147
$ws.Range("P3").Value = 0.7588727250690147
$ws.Range("Q3").Value = 0.6079098977575556
$ws.Range("R3").Value = 5.471189079818
$ws.Range("S3").Value = 0.1027395858550824
$ws.Range("T3").Value = 0.1027395858550824

# Row 4
$ws.Range("I4").Value = 0.389303862711544
$ws.Range("J4").Value = 0.389303862711544
$ws.Range("M4").Value = 0.4655266666666667
$ws.Range("N4").Value = 1.39658
$ws.Range("O4").Value = 0.2411272749309853
$ws.Range("P4").Value = 0.2411272749309853
$ws.Range("Q4").Value = 0.5554391077688889
$ws.Range("R4").Value = 4.99895196992
$ws.Range("S4").Value = 0.09387177953574105
$ws.Range("T4").Value = 0.09387177953574104

# Row 5
$ws.Range("I5").Value = 0.389303862711544
$ws.Range("J5").Value = 0.389303862711544
$ws.Range("M5").Value = 1.465099666666666
$ws.Range("N5").Value = 4.395299
$ws.Range("O5").Value = 0.7588727250690147
$ws.Range("P5").Value = 0.7588727250690147
$ws.Range("Q5").Value = 1.748070969752889
$ws.Range("R5").Value = 15.732638727776
$ws.Range("S5").Value = 0.295432083175803
$ws.Range("T5").Value = 0.295432083175803

# Row 6
$ws.Range("G6").Value = 1.270157666666667
$ws.Range("H6").Value = 3.810473
$ws.Range("I6").Value = 0.4144331204288861
$ws.Range("J6").Value = 0.4144331204288861
$ws.Range("M6").Value = 0.4655266666666667
$ws.Range("N6").Value = 1.39658
$ws.Range("O6").Value = 0.2411272749309853
$ws.Range("P6").Value = 0.2411272749309853
$ws.Range("Q6").Value = 0.5912922647044445
$ws.Range("R6").Value = 5.32163038234
$ws.Range("S6").Value = 0.09993112897016218
$ws.Range("T6").Value = 0.09993112897016217

# Row 7
$ws.Range("G7").Value = 1.270157666666667
$ws.Range("H7").Value = 3.810473
$ws.Range("I7").Value = 0.4144331204288861
$ws.Range("J7").Value = 0.4144331204288861
$ws.Range("M7").Value = 1.465099666666666
$ws.Range("N7").Value = 4.395299
$ws.Range("O7").Value = 0.7588727250690147
$ws.Range("P7").Value = 0.7588727250690147
$ws.Range("Q7").Value = 1.860907574047444
$ws.Range("R7").Value = 16.748168166427
$ws.Range("S7").Value = 0.314501991458724
$ws.Range("T7").Value = 0.314501991458724

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.186581
$ws.Range("H8").Value = 0.559743
$ws.Range("I8").Value = 0.06087854135909794
$ws.Range("J8").Value = 0.06087854135909794
$ws.Range("M8").Value = 0.4655266666666667
$ws.Range("N8").Value = 1.39658
$ws.Range("O8").Value = 0.2411272749309853
$ws.Range("P8").Value = 0.2411272749309853
$ws.Range("Q8").Value = 0.08685843099333333
$ws.Range("R8").Value = 0.7817258789400001
$ws.Range("S8").Value = 0.01467947677969257
$ws.Range("T8").Value = 0.01467947677969257

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.186581
$ws.Range("H9").Value = 0.559743
$ws.Range("I9").Value = 0.06087854135909794
$ws.Range("J9").Value = 0.06087854135909794
$ws.Range("M9").Value = 1.465099666666666
$ws.Range("N9").Value = 4.395299
$ws.Range("O9").Value = 0.7588727250690147
$ws.Range("P9").Value = 0.7588727250690147
$ws.Range("Q9").Value = 0.2733597609063333
$ws.Range("R9").Value = 2.460237848157
$ws.Range("S9").Value = 0.04619906457940537
$ws.Range("T9").Value = 0.04619906457940537
